$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 67, shifting existing rows 67-78 down to 68-79.
$ws.Rows.Item(67).Insert()

# Fill in the newly inserted row 67 with the new weekly record.
$ws.Range("A67").Value = 10
$ws.Range("B67").Value = "Vega Modelo de Temuco"
$ws.Range("C67").Value = "La Araucanía"
$ws.Range("D67").Value = 44551
$ws.Range("D67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E67").Value = 9
$ws.Range("F67").Value = "Fruta"
$ws.Range("G67").Value = 100101
$ws.Range("H67").Value = "Berries"
$ws.Range("I67").Value = 100101001
$ws.Range("J67").Value = "Arándano (blue)"
$ws.Range("K67").Value = "Sin especificar"
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 300
$ws.Range("N67").Value = 2600
$ws.Range("O67").Value = 2600
$ws.Range("P67").Value = 2600
$ws.Range("Q67").Value = "$/kilo"
$ws.Range("R67").Value = "Región de O'Higgins"
$ws.Range("S67").Value = 2600
$ws.Range("T67").Value = 1
